# Refresh the crypto price/volume snapshot in the "cryptos" sheet.
# Mirrors the upstream GitHub Actions data pull: most rows just get new
# Price (D) / Volume(1h) (E) quotes, a few rows swap which coin occupies
# that rank (B/C/D/E all change) while the rank index in column A is left
# untouched, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '34.629.80'
$ws.Cells.Item(2, 5).Value = '  +2.63%  '

# Row 3: Ethereum
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.789.34'
$ws.Cells.Item(3, 5).Value = '  +0.90%  '

# Row 4: TetherUSD
$ws.Cells.Item(4, 5).Value = '  -0.22%  '

# Row 5: BNB
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '224.42'
$ws.Cells.Item(5, 5).Value = '  -0.12%  '

# Row 6: XRP
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.558'
$ws.Cells.Item(6, 5).Value = '  +0.92%  '

# Row 7: USDC
$ws.Cells.Item(7, 5).Value = '  -0.11%  '

# Row 8: Solana
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '32.97'
$ws.Cells.Item(8, 5).Value = '  +7.87%  '

# Row 9: Cardano
$ws.Cells.Item(9, 5).Value = '  +2.26%  '

# Row 10: Dogecoin
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.0678'
$ws.Cells.Item(10, 5).Value = '  +3.07%  '

# Row 11: TRON
$ws.Cells.Item(11, 5).Value = '  +1.31%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '2.048.51'
$ws.Cells.Item(12, 5).Value = '  +0.96%  '

# Row 13: Chainlink
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '11.02'
$ws.Cells.Item(13, 5).Value = '  +11.03%  '

# Row 14: WrappedEther
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '1.784.59'
$ws.Cells.Item(14, 5).Value = '  +0.67%  '

# Row 15: Polygon
$ws.Cells.Item(15, 2).Value = 'Polygon'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.634'
$ws.Cells.Item(15, 5).Value = '  +1.73%  '

# Row 16: WrappedBTC
$ws.Cells.Item(16, 2).Value = 'WrappedBTC'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '34.620.38'
$ws.Cells.Item(16, 5).Value = '  +2.64%  '

# Row 17: Polkadot
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '4.29'
$ws.Cells.Item(17, 5).Value = '  +2.84%  '

# Row 18: Litecoin
$ws.Cells.Item(18, 5).Value = '  +0.62%  '

# Row 19: BitcoinCash
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '253.85'
$ws.Cells.Item(19, 5).Value = '  +1.17%  '

# Row 20: ShibaInu
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '0.0₃0772'
$ws.Cells.Item(20, 5).Value = '  +4.93%  '

# Row 21: Dai
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '0.999'
$ws.Cells.Item(21, 5).Value = '  -0.42%  '

# Row 22: Avalanche
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '10.41'
$ws.Cells.Item(22, 5).Value = '  +1.85%  '

# Row 23: Uniswap
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '4.24'
$ws.Cells.Item(23, 5).Value = '  +1.46%  '

# Row 24: Toncoin
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '2.12'
$ws.Cells.Item(24, 5).Value = '  -1.05%  '

# Row 25: Monero
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '159.29'
$ws.Cells.Item(25, 5).Value = '  +0.51%  '

# Row 26: EthereumClassic
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '16.37'
$ws.Cells.Item(26, 5).Value = '  -0.15%  '

# Row 27: Cosmos
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '7.10'
$ws.Cells.Item(27, 5).Value = '  +3.09%  '

# Row 28: Stellar
$ws.Cells.Item(28, 5).Value = '  +0.59%  '

# Row 29: BinanceUSD
$ws.Cells.Item(29, 5).Value = '  -0.20%  '

# Row 30: Filecoin
$ws.Cells.Item(30, 2).Value = 'Filecoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '3.76'
$ws.Cells.Item(30, 5).Value = '  -0.64%  '

# Row 31: Hedera
$ws.Cells.Item(31, 2).Value = 'Hedera'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.0516'
$ws.Cells.Item(31, 5).Value = '  +1.03%  '

# Row 32: PancakeSwap
$ws.Cells.Item(32, 5).Value = '  +0.48%  '

# Row 33: InternetComputer(DFINITY)
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '3.59'
$ws.Cells.Item(33, 5).Value = '  +1.64%  '

# Row 34: LidoDAOToken
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.86'
$ws.Cells.Item(34, 5).Value = '  +4.25%  '

# Row 35: Maker
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.443.49'
$ws.Cells.Item(35, 5).Value = '  -2.17%  '

# Row 36: TrustWalletToken
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '1.05'
$ws.Cells.Item(36, 5).Value = '  -0.40%  '

# Row 37: VeChain
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.0190'
$ws.Cells.Item(37, 5).Value = '  +2.98%  '

# Row 38: ImmutableX
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.629'
$ws.Cells.Item(38, 5).Value = '  -0.04%  '

# Row 39: Aave
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '83.27'
$ws.Cells.Item(39, 5).Value = '  +0.28%  '

# Row 40: MXToken
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '2.81'
$ws.Cells.Item(40, 5).Value = '  +4.68%  '

# Row 41: HuobiToken
$ws.Cells.Item(41, 5).Value = '  -0.29%  '

# Row 42: ARBITRUM
$ws.Cells.Item(42, 5).Value = '  +1.90%  '

# Row 43: RenderToken
$ws.Cells.Item(43, 5).Value = '  +0.04%  '

# Row 44: Kaspa
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.0501'
$ws.Cells.Item(44, 5).Value = '  -1.06%  '

# Row 45: FraxShare
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '5.90'
$ws.Cells.Item(45, 5).Value = '  +2.76%  '

# Row 46: WEMIXToken
$ws.Cells.Item(46, 5).Value = '  -1.74%  '

# Row 47: RocketPoolETH
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.944.46'
$ws.Cells.Item(47, 5).Value = '  +0.65%  '

# Row 48: Quant
$ws.Cells.Item(48, 2).Value = 'Quant'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '104.15'
$ws.Cells.Item(48, 5).Value = '  +7.16%  '

# Row 49: InjectiveProtocol
$ws.Cells.Item(49, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '11.99'
$ws.Cells.Item(49, 5).Value = '  +0.72%  '

# Row 50: PaxDollar
$ws.Cells.Item(50, 2).Value = 'PaxDollar'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '1.00'
$ws.Cells.Item(50, 5).Value = '  -0.20%  '

# Row 51: BabyDogeCoin
$ws.Cells.Item(51, 5).Value = '  +5.39%  '
